$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "author"
$ws.Range("B1").Value = "titulo"

$ws.Range("A2").Value = 1

$ws.Range("B1").Select()
